$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Cells.Item(7, 1).Value = "a1"
$ws.Cells.Item(7, 2).Value = 0
$ws.Cells.Item(7, 3).Value = 2.598305842476279
$ws.Cells.Item(7, 4).Value = 51.50543951451338
$ws.Cells.Item(7, 5).Value = 31.16742599597234
$ws.Cells.Item(7, 6).Value = 4.258863675201241
$ws.Cells.Item(7, 7).Value = 1.056793505147231
$ws.Cells.Item(7, 8).Value = 0.06353140626154063
$ws.Cells.Item(7, 9).Value = 9.320671523399811
$ws.Cells.Item(7, 10).Value = 0.01023373335740843
$ws.Cells.Item(7, 11).Value = 0.004158830126264144
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = 0.0001285637541733014
$ws.Cells.Item(7, 14).Value = 0.01444400804360804
$ws.Cells.Item(7, 15).Value = [double]"3.323798869188713e-06"
$ws.Cells.Item(7, 16).Value = [double]"4.353728220697508e-09"
$ws.Cells.Item(7, 17).Value = [double]"1.961925133339872e-09"
$ws.Cells.Item(7, 18).Value = [double]"8.485590254492418e-09"
$ws.Cells.Item(7, 19).Value = [double]"2.517200017541414e-08"
$ws.Cells.Item(7, 20).Value = [double]"3.642096472706921e-08"
$ws.Cells.Item(7, 21).Value = [double]"1.501676190244253e-09"
$ws.Cells.Item(7, 22).Value = [double]"4.807390552682693e-11"
$ws.Cells.Item(7, 23).Value = [double]"1.840474756842591e-13"
$ws.Cells.Item(7, 24).Value = [double]"4.078040261658826e-14"
$ws.Cells.Item(7, 25).Value = [double]"9.995168382840823e-15"
$ws.Cells.Item(7, 26).Value = [double]"2.408840620969784e-12"
$ws.Cells.Item(7, 27).Value = [double]"7.541791353791227e-13"
$ws.Cells.Item(7, 28).Value = [double]"1.408989222118458e-13"
$ws.Cells.Item(7, 29).Value = [double]"3.206723940839607e-13"
$ws.Cells.Item(7, 30).Value = [double]"3.996604204341381e-14"
$ws.Cells.Item(7, 31).Value = 0
$ws.Cells.Item(7, 32).Value = [double]"2.582391302583239e-16"
$ws.Cells.Item(7, 33).Value = [double]"2.565608363901613e-17"
$ws.Cells.Item(7, 34).Value = [double]"2.308540254639245e-20"
$ws.Cells.Item(7, 35).Value = [double]"6.013616547206893e-19"
$ws.Cells.Item(7, 36).Value = [double]"1.20151189943135e-17"
$ws.Cells.Item(7, 37).Value = 0
$ws.Cells.Item(7, 38).Value = [double]"6.724142812192292e-21"
$ws.Cells.Item(7, 39).Value = [double]"1.921807528063854e-22"
$ws.Cells.Item(7, 40).Value = [double]"4.469866205371541e-27"
$ws.Cells.Item(7, 41).Value = 0
$ws.Cells.Item(7, 42).Value = 0
$ws.Cells.Item(7, 43).Value = 0
$ws.Cells.Item(7, 44).Value = 0.01454847919631769
$ws.Cells.Item(7, 45).Value = 1.598023128028724
$ws.Cells.Item(7, 46).Value = 3.294500958061927
$ws.Cells.Item(7, 47).Value = 1.597484233197451
$ws.Cells.Item(7, 48).Value = 1.004584623398545
$ws.Cells.Item(7, 49).Value = 0.1936103460199946
$ws.Cells.Item(7, 50).Value = 90.2294905325778
$ws.Cells.Item(7, 51).Value = 0.06184189372697138
$ws.Cells.Item(7, 52).Value = 0.05608128015492498
$ws.Cells.Item(7, 53).Value = 0
$ws.Cells.Item(7, 54).Value = 0.005234118075935938
$ws.Cells.Item(7, 55).Value = 1.944220010668491
$ws.Cells.Item(7, 56).Value = 0.0003033541557622196
$ws.Cells.Item(7, 57).Value = [double]"1.75434494150938e-06"
$ws.Cells.Item(7, 58).Value = [double]"8.299458692300326e-07"
$ws.Cells.Item(7, 59).Value = [double]"7.40282051373091e-06"
$ws.Cells.Item(7, 60).Value = [double]"2.541722905693818e-05"
$ws.Cells.Item(7, 61).Value = [double]"3.943184041252878e-05"
$ws.Cells.Item(7, 62).Value = [double]"2.143205035036791e-06"
$ws.Cells.Item(7, 63).Value = [double]"5.43888709694314e-08"
$ws.Cells.Item(7, 64).Value = [double]"3.412872884437011e-10"
$ws.Cells.Item(7, 65).Value = [double]"4.803981711680048e-11"
$ws.Cells.Item(7, 66).Value = [double]"2.804378191395179e-11"
$ws.Cells.Item(7, 67).Value = [double]"5.272136487519178e-09"
$ws.Cells.Item(7, 68).Value = [double]"1.721122076908615e-09"
$ws.Cells.Item(7, 69).Value = [double]"6.428602440131308e-10"
$ws.Cells.Item(7, 70).Value = [double]"6.89240488159044e-10"
$ws.Cells.Item(7, 71).Value = [double]"2.174220148391036e-10"
$ws.Cells.Item(7, 72).Value = 0
$ws.Cells.Item(7, 73).Value = [double]"2.025844747059195e-12"
$ws.Cells.Item(7, 74).Value = [double]"1.392140110576935e-13"
$ws.Cells.Item(7, 75).Value = [double]"2.159498603293262e-16"
$ws.Cells.Item(7, 76).Value = [double]"4.420182735185514e-15"
$ws.Cells.Item(7, 77).Value = [double]"1.264756660063012e-13"
$ws.Cells.Item(7, 78).Value = 0
$ws.Cells.Item(7, 79).Value = [double]"1.802184769970499e-16"
$ws.Cells.Item(7, 80).Value = [double]"6.260858938615566e-18"
$ws.Cells.Item(7, 81).Value = [double]"4.528437939613383e-22"
$ws.Cells.Item(7, 82).Value = 0
$ws.Cells.Item(7, 83).Value = 0
$ws.Cells.Item(7, 84).Value = 0
$ws.Cells.Item(7, 85).Value = 0.0002711975470577192
$ws.Cells.Item(7, 86).Value = 0.01762952303297098
$ws.Cells.Item(7, 87).Value = 0.03404198362790177
$ws.Cells.Item(7, 88).Value = 0.02397282665632248
$ws.Cells.Item(7, 89).Value = 0.02756530111955153
$ws.Cells.Item(7, 90).Value = 0.0201883247375805
$ws.Cells.Item(7, 91).Value = 79.73618264289385
$ws.Cells.Item(7, 92).Value = 0.02729720732264617
$ws.Cells.Item(7, 93).Value = 0.06441826959471283
$ws.Cells.Item(7, 94).Value = 0
$ws.Cells.Item(7, 95).Value = 0.02518353031595424
$ws.Cells.Item(7, 96).Value = 20.01767083148066
$ws.Cells.Item(7, 97).Value = 0.001788052637497496
$ws.Cells.Item(7, 98).Value = [double]"2.950002153546877e-05"
$ws.Cells.Item(7, 99).Value = [double]"1.678144141370301e-05"
$ws.Cells.Item(7, 100).Value = 0.0002865976628999058
$ws.Cells.Item(7, 101).Value = 0.001204097416335504
$ws.Cells.Item(7, 102).Value = 0.002042934613753144
$ws.Cells.Item(7, 103).Value = 0.0002039208298505403
$ws.Cells.Item(7, 104).Value = [double]"2.701463549556472e-06"
$ws.Cells.Item(7, 105).Value = [double]"7.158817026616794e-08"
$ws.Cells.Item(7, 106).Value = [double]"1.075680471894077e-08"
$ws.Cells.Item(7, 107).Value = [double]"1.299405880782798e-08"
$ws.Cells.Item(7, 108).Value = [double]"1.793177415731068e-06"
$ws.Cells.Item(7, 109).Value = [double]"6.69086853382651e-07"
$ws.Cells.Item(7, 110).Value = [double]"5.680119241240874e-07"
$ws.Cells.Item(7, 111).Value = [double]"2.938730439796412e-07"
$ws.Cells.Item(7, 112).Value = [double]"3.470911221708254e-07"
$ws.Cells.Item(7, 113).Value = 0
$ws.Cells.Item(7, 114).Value = [double]"7.56941570895172e-09"
$ws.Cells.Item(7, 115).Value = [double]"4.465596245744859e-10"
$ws.Cells.Item(7, 116).Value = [double]"1.434894372053936e-12"
$ws.Cells.Item(7, 117).Value = [double]"2.982133273447261e-11"
$ws.Cells.Item(7, 118).Value = [double]"9.465916868629856e-10"
$ws.Cells.Item(7, 119).Value = 0
$ws.Cells.Item(7, 120).Value = [double]"1.007430280322403e-11"
$ws.Cells.Item(7, 121).Value = [double]"6.930788694037974e-13"
$ws.Cells.Item(7, 122).Value = [double]"6.836784561793773e-16"
$ws.Cells.Item(7, 123).Value = 0
$ws.Cells.Item(7, 124).Value = 0
$ws.Cells.Item(7, 125).Value = 0
$ws.Cells.Item(7, 126).Value = 0
$ws.Cells.Item(7, 127).Value = 0
$ws.Cells.Item(7, 128).Value = 0
$ws.Cells.Item(7, 129).Value = [double]"3.418568208555885e-24"
$ws.Cells.Item(7, 130).Value = [double]"1.820345185165017e-16"
$ws.Cells.Item(7, 131).Value = [double]"1.593090825193606e-09"
$ws.Cells.Item(7, 132).Value = 0.2382122479497648
$ws.Cells.Item(7, 133).Value = 0.0001798624602656845
$ws.Cells.Item(7, 134).Value = 0.00369814131969511
$ws.Cells.Item(7, 135).Value = 0
$ws.Cells.Item(7, 136).Value = 0.01069063773472701
$ws.Cells.Item(7, 137).Value = 54.25741044704761
$ws.Cells.Item(7, 138).Value = 0.009920251244325138
$ws.Cells.Item(7, 139).Value = 0.01205680495179951
$ws.Cells.Item(7, 140).Value = 0.01205903728671882
$ws.Cells.Item(7, 141).Value = 0.9202506795587261
$ws.Cells.Item(7, 142).Value = 6.320651755777595
$ws.Cells.Item(7, 143).Value = 13.3902997863129
$ws.Cells.Item(7, 144).Value = 5.932368369783462
$ws.Cells.Item(7, 145).Value = 0.01750537921359408
$ws.Cells.Item(7, 146).Value = 0.01491989534584642
$ws.Cells.Item(7, 147).Value = 0.003309029307115738
$ws.Cells.Item(7, 148).Value = 0.01836296222646635
$ws.Cells.Item(7, 149).Value = 1.239499929593332
$ws.Cells.Item(7, 150).Value = 0.6495897091817214
$ws.Cells.Item(7, 151).Value = 3.49699655389669
$ws.Cells.Item(7, 152).Value = 0.3649638582949725
$ws.Cells.Item(7, 153).Value = 8.931286484695583
$ws.Cells.Item(7, 154).Value = 0
$ws.Cells.Item(7, 155).Value = 1.425425039019363
$ws.Cells.Item(7, 156).Value = 0.0642703772104268
$ws.Cells.Item(7, 157).Value = 0.001147685275968893
$ws.Cells.Item(7, 158).Value = 0.02524907358759119
$ws.Cells.Item(7, 159).Value = 0.9583171288525215
$ws.Cells.Item(7, 160).Value = 0
$ws.Cells.Item(7, 161).Value = 1.104073195366162
$ws.Cells.Item(7, 162).Value = 0.4039852147827664
$ws.Cells.Item(7, 163).Value = 0.1733004611292089
$ws.Cells.Item(7, 164).Value = 0
$ws.Cells.Item(7, 165).Value = 0

# Row 8
$ws.Cells.Item(8, 1).Value = "b2"
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(8, 3).Value = 2.598307891396049
$ws.Cells.Item(8, 4).Value = 51.50548012963711
$ws.Cells.Item(8, 5).Value = 31.16745057335436
$ws.Cells.Item(8, 6).Value = 4.258867033570313
$ws.Cells.Item(8, 7).Value = 1.056794338492202
$ws.Cells.Item(8, 8).Value = 0.06353145635985938
$ws.Cells.Item(8, 9).Value = 9.320678873307525
$ws.Cells.Item(8, 10).Value = 0.01023374569252173
$ws.Cells.Item(8, 11).Value = 0.004158833405750771
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = 0.0001285081147045908
$ws.Cells.Item(8, 14).Value = 0.01436526848407125
$ws.Cells.Item(8, 15).Value = [double]"3.270279599084137e-06"
$ws.Cells.Item(8, 16).Value = [double]"4.346703944849269e-09"
$ws.Cells.Item(8, 17).Value = [double]"1.959679230772089e-09"
$ws.Cells.Item(8, 18).Value = [double]"8.481474297833965e-09"
$ws.Cells.Item(8, 19).Value = [double]"2.516294048446699e-08"
$ws.Cells.Item(8, 20).Value = [double]"3.640203817462737e-08"
$ws.Cells.Item(8, 21).Value = [double]"1.50132966879395e-09"
$ws.Cells.Item(8, 22).Value = [double]"4.787413534413458e-11"
$ws.Cells.Item(8, 23).Value = [double]"1.839793747843049e-13"
$ws.Cells.Item(8, 24).Value = [double]"4.072635764951946e-14"
$ws.Cells.Item(8, 25).Value = [double]"9.991437339627481e-15"
$ws.Cells.Item(8, 26).Value = [double]"2.407997174901375e-12"
$ws.Cells.Item(8, 27).Value = [double]"7.538016255571783e-13"
$ws.Cells.Item(8, 28).Value = [double]"1.408516410435552e-13"
$ws.Cells.Item(8, 29).Value = [double]"3.204624978744216e-13"
$ws.Cells.Item(8, 30).Value = [double]"3.995191568085481e-14"
$ws.Cells.Item(8, 31).Value = 0
$ws.Cells.Item(8, 32).Value = [double]"2.581990920774012e-16"
$ws.Cells.Item(8, 33).Value = [double]"2.565276326216815e-17"
$ws.Cells.Item(8, 34).Value = [double]"2.307626514970993e-20"
$ws.Cells.Item(8, 35).Value = [double]"6.011460822270416e-19"
$ws.Cells.Item(8, 36).Value = [double]"1.20089892050355e-17"
$ws.Cells.Item(8, 37).Value = 0
$ws.Cells.Item(8, 38).Value = [double]"6.722847296857941e-21"
$ws.Cells.Item(8, 39).Value = [double]"1.921023791287214e-22"
$ws.Cells.Item(8, 40).Value = [double]"4.466358367442952e-27"
$ws.Cells.Item(8, 41).Value = 0
$ws.Cells.Item(8, 42).Value = 0
$ws.Cells.Item(8, 43).Value = 0
$ws.Cells.Item(8, 44).Value = 0.01454849887653683
$ws.Cells.Item(8, 45).Value = 1.598025289728578
$ws.Cells.Item(8, 46).Value = 3.294505414644637
$ws.Cells.Item(8, 47).Value = 1.597486394168324
$ws.Cells.Item(8, 48).Value = 1.004585982334093
$ws.Cells.Item(8, 49).Value = 0.1936106079232483
$ws.Cells.Item(8, 50).Value = 90.22961258905691
$ws.Cells.Item(8, 51).Value = 0.0618420037509636
$ws.Cells.Item(8, 52).Value = 0.05608135601796665
$ws.Cells.Item(8, 53).Value = 0
$ws.Cells.Item(8, 54).Value = 0.00523590464469739
$ws.Cells.Item(8, 55).Value = 1.944085575763411
$ws.Cells.Item(8, 56).Value = 0.0003033402518602086
$ws.Cells.Item(8, 57).Value = [double]"1.754346592996569e-06"
$ws.Cells.Item(8, 58).Value = [double]"8.299469197583931e-07"
$ws.Cells.Item(8, 59).Value = [double]"7.40282915122265e-06"
$ws.Cells.Item(8, 60).Value = [double]"2.541726343968232e-05"
$ws.Cells.Item(8, 61).Value = [double]"3.943189375331102e-05"
$ws.Cells.Item(8, 62).Value = [double]"2.143207518277915e-06"
$ws.Cells.Item(8, 63).Value = [double]"5.438892817768887e-08"
$ws.Cells.Item(8, 64).Value = [double]"3.412877501145423e-10"
$ws.Cells.Item(8, 65).Value = [double]"4.803987726835131e-11"
$ws.Cells.Item(8, 66).Value = [double]"2.804381817328833e-11"
$ws.Cells.Item(8, 67).Value = [double]"5.272142261404365e-09"
$ws.Cells.Item(8, 68).Value = [double]"1.721124405128575e-09"
$ws.Cells.Item(8, 69).Value = [double]"6.428610737249032e-10"
$ws.Cells.Item(8, 70).Value = [double]"6.892413538796666e-10"
$ws.Cells.Item(8, 71).Value = [double]"2.174223089532059e-10"
$ws.Cells.Item(8, 72).Value = 0
$ws.Cells.Item(8, 73).Value = [double]"2.025847487487802e-12"
$ws.Cells.Item(8, 74).Value = [double]"1.392141993771879e-13"
$ws.Cells.Item(8, 75).Value = [double]"2.159501524519954e-16"
$ws.Cells.Item(8, 76).Value = [double]"4.420188714515979e-15"
$ws.Cells.Item(8, 77).Value = [double]"1.26475837094206e-13"
$ws.Cells.Item(8, 78).Value = 0
$ws.Cells.Item(8, 79).Value = [double]"1.802186967531235e-16"
$ws.Cells.Item(8, 80).Value = [double]"6.260867407890897e-18"
$ws.Cells.Item(8, 81).Value = [double]"4.52844406538432e-22"
$ws.Cells.Item(8, 82).Value = 0
$ws.Cells.Item(8, 83).Value = 0
$ws.Cells.Item(8, 84).Value = 0
$ws.Cells.Item(8, 85).Value = 0.0002714386214454296
$ws.Cells.Item(8, 86).Value = 0.01764519436376623
$ws.Cells.Item(8, 87).Value = 0.03407224441177899
$ws.Cells.Item(8, 88).Value = 0.02399413671082172
$ws.Cells.Item(8, 89).Value = 0.02758980461584619
$ws.Cells.Item(8, 90).Value = 0.02020627065220184
$ws.Cells.Item(8, 91).Value = 79.80706216085981
$ws.Cells.Item(8, 92).Value = 0.02732148430247454
$ws.Cells.Item(8, 93).Value = 0.06447553263071149
$ws.Cells.Item(8, 94).Value = 0
$ws.Cells.Item(8, 95).Value = 0.02518872388660588
$ws.Cells.Item(8, 96).Value = 19.94659585337205
$ws.Cells.Item(8, 97).Value = 0.001784723507098906
$ws.Cells.Item(8, 98).Value = [double]"2.950714429183297e-05"
$ws.Cells.Item(8, 99).Value = [double]"1.678926177057798e-05"
$ws.Cells.Item(8, 100).Value = 0.0002867629402797082
$ws.Cells.Item(8, 101).Value = 0.001204716553711497
$ws.Cells.Item(8, 102).Value = 0.002044151224989888
$ws.Cells.Item(8, 103).Value = 0.0002040248329032531
$ws.Cells.Item(8, 104).Value = [double]"2.702621103372722e-06"
$ws.Cells.Item(8, 105).Value = [double]"7.161593733445081e-08"
$ws.Cells.Item(8, 106).Value = [double]"1.076384412296298e-08"
$ws.Cells.Item(8, 107).Value = [double]"1.300188557112609e-08"
$ws.Cells.Item(8, 108).Value = [double]"1.793959030106487e-06"
$ws.Cells.Item(8, 109).Value = [double]"6.694598335151108e-07"
$ws.Cells.Item(8, 110).Value = [double]"5.683184096622083e-07"
$ws.Cells.Item(8, 111).Value = [double]"2.940598335068914e-07"
$ws.Cells.Item(8, 112).Value = [double]"3.472980855063064e-07"
$ws.Cells.Item(8, 113).Value = 0
$ws.Cells.Item(8, 114).Value = [double]"7.573674935699098e-09"
$ws.Cells.Item(8, 115).Value = [double]"4.468338769863724e-10"
$ws.Cells.Item(8, 116).Value = [double]"1.435629875215994e-12"
$ws.Cells.Item(8, 117).Value = [double]"2.982822350000825e-11"
$ws.Cells.Item(8, 118).Value = [double]"9.469325872712386e-10"
$ws.Cells.Item(8, 119).Value = 0
$ws.Cells.Item(8, 120).Value = [double]"1.007894729562178e-11"
$ws.Cells.Item(8, 121).Value = [double]"6.93404095953161e-13"
$ws.Cells.Item(8, 122).Value = [double]"6.835895990840956e-16"
$ws.Cells.Item(8, 123).Value = 0
$ws.Cells.Item(8, 124).Value = 0
$ws.Cells.Item(8, 125).Value = 0
$ws.Cells.Item(8, 126).Value = 0
$ws.Cells.Item(8, 127).Value = 0
$ws.Cells.Item(8, 128).Value = 0
$ws.Cells.Item(8, 129).Value = [double]"3.365436763557885e-24"
$ws.Cells.Item(8, 130).Value = [double]"1.792628436493903e-16"
$ws.Cells.Item(8, 131).Value = [double]"1.575252139961857e-09"
$ws.Cells.Item(8, 132).Value = 0.2381758563270786
$ws.Cells.Item(8, 133).Value = 0.0001798310664168813
$ws.Cells.Item(8, 134).Value = 0.003697582366629099
$ws.Cells.Item(8, 135).Value = 0
$ws.Cells.Item(8, 136).Value = 0.01069180841002004
$ws.Cells.Item(8, 137).Value = 54.26431735892116
$ws.Cells.Item(8, 138).Value = 0.009919589875755887
$ws.Cells.Item(8, 139).Value = 0.01205498714872031
$ws.Cells.Item(8, 140).Value = 0.01205721627104594
$ws.Cells.Item(8, 141).Value = 0.9201117688308835
$ws.Cells.Item(8, 142).Value = 6.319696907627401
$ws.Cells.Item(8, 143).Value = 13.38827694065148
$ws.Cells.Item(8, 144).Value = 5.93147217882357
$ws.Cells.Item(8, 145).Value = 0.01750273679770274
$ws.Cells.Item(8, 146).Value = 0.01491764142725977
$ws.Cells.Item(8, 147).Value = 0.003308529418712624
$ws.Cells.Item(8, 148).Value = 0.0183601881706899
$ws.Cells.Item(8, 149).Value = 1.239312680831548
$ws.Cells.Item(8, 150).Value = 0.6494915769706515
$ws.Cells.Item(8, 151).Value = 3.496468269659598
$ws.Cells.Item(8, 152).Value = 0.3649087239388261
$ws.Cells.Item(8, 153).Value = 8.929937253206763
$ws.Cells.Item(8, 154).Value = 0
$ws.Cells.Item(8, 155).Value = 1.425209702925185
$ws.Cells.Item(8, 156).Value = 0.06426066801378651
$ws.Cells.Item(8, 157).Value = 0.001147511897462197
$ws.Cells.Item(8, 158).Value = 0.02524525926081909
$ws.Cells.Item(8, 159).Value = 0.9581723578109992
$ws.Cells.Item(8, 160).Value = 0
$ws.Cells.Item(8, 161).Value = 1.10390640524878
$ws.Cells.Item(8, 162).Value = 0.4039241855487656
$ws.Cells.Item(8, 163).Value = 0.1732742809770457
$ws.Cells.Item(8, 164).Value = 0
$ws.Cells.Item(8, 165).Value = 0

# Match column-A label style (bold, centered, bordered) copied from the existing A6 label cell
$ws.Cells.Item(6, 1).Copy()
$ws.Cells.Item(7, 1).PasteSpecial(-4122)
$ws.Cells.Item(6, 1).Copy()
$ws.Cells.Item(8, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
